$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 214
$ws1.Range("F4").Value = 398
$ws1.Range("F5").Value = 196
$ws1.Range("F6").Value = 785
$ws1.Range("F8").Value = 10101
$ws1.Range("F10").Value = 3480
$ws1.Range("F12").Value = 2429
$ws1.Range("F14").Value = 2771
$ws1.Range("F17").Value = 2144
$ws1.Range("F22").Value = 22
$ws1.Range("F24").Value = 307
$ws1.Range("F26").Value = 207
$ws1.Range("F28").Value = 1310
$ws1.Range("F29").Value = 6
$ws1.Range("F34").Value = 2967
$ws1.Range("F35").Value = 2954
$ws1.Range("F36").Value = 20
$ws1.Range("F38").Value = 1027
$ws1.Range("F39").Value = 387
$ws1.Range("F41").Value = 1289
$ws1.Range("F43").Value = 103
$ws1.Range("F44").Value = 69

# Sheet "演出" (Performance) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 176
$ws2.Range("F16").Value = 173

# Sheet "本地生活" (Local Life) - column F updates
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 740
$ws3.Range("F3").Value = 979
$ws3.Range("F5").Value = 1978

# Sheet "全部类型" (All Types) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 979
$ws4.Range("F5").Value = 398
$ws4.Range("F8").Value = 196
$ws4.Range("F10").Value = 10101
$ws4.Range("F12").Value = 3480
$ws4.Range("F17").Value = 2144
$ws4.Range("F22").Value = 307
$ws4.Range("F24").Value = 207
$ws4.Range("F26").Value = 1310
$ws4.Range("F27").Value = 6
$ws4.Range("F33").Value = 2967
$ws4.Range("F34").Value = 2954
$ws4.Range("F35").Value = 20
$ws4.Range("F36").Value = 1027
$ws4.Range("F39").Value = 387
$ws4.Range("F45").Value = 69
$ws4.Range("F49").Value = 173
